# Insert a new data row at row 62 (pushing existing rows 62..148 down to 63..149)
# and populate the newly inserted row with the latest weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 62; this shifts formatting/content of
# rows 62..148 down to 63..149, matching the row immediately above (row 61)
# for formatting purposes (same as every other data row in this table).
$ws.Rows(62).Insert()

# Populate the new row 62 with the new weekly observation.
$ws.Cells.Item(62, 1).Value  = 4
$ws.Cells.Item(62, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(62, 3).Value  = "Los Lagos"
$ws.Cells.Item(62, 4).Value  = 44495
$ws.Cells.Item(62, 5).Value  = 10
$ws.Cells.Item(62, 6).Value  = "Fruta"
$ws.Cells.Item(62, 7).Value  = 100108
$ws.Cells.Item(62, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(62, 9).Value  = 100108005
$ws.Cells.Item(62, 10).Value = "Piña"
$ws.Cells.Item(62, 11).Value = "Caramelo"
$ws.Cells.Item(62, 12).Value = "Segunda"
$ws.Cells.Item(62, 13).Value = 270
$ws.Cells.Item(62, 14).Value = 23500
$ws.Cells.Item(62, 15).Value = 24000
$ws.Cells.Item(62, 16).Value = 23778
$ws.Cells.Item(62, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(62, 18).Value = "Ecuador"
$ws.Cells.Item(62, 19).Value = 1698
$ws.Cells.Item(62, 20).Value = 14

# Ensure the date cell keeps the same date-time number format used by the
# rest of the "Fecha" column (style index 2 in the original workbook).
$ws.Cells.Item(62, 4).NumberFormat = $ws.Cells.Item(63, 4).NumberFormat
